# Optuna Attempt (go back with original)
# Apply updated forecast values to "Forecast Comparison" and "Summary" sheets.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2
$wsForecast.Range("D2").Value = 86
$wsForecast.Range("H2").Value = 14.93
$wsForecast.Range("L2").Value = 1.08

# Row 3
$wsForecast.Range("D3").Value = 123
$wsForecast.Range("H3").Value = 9.699999999999999
$wsForecast.Range("L3").Value = 0.85

# Row 4
$wsForecast.Range("H4").Value = 11.81
$wsForecast.Range("L4").Value = 0.86

# Row 5
$wsForecast.Range("H5").Value = 9.460000000000001
$wsForecast.Range("L5").Value = 0.92

# Row 6
$wsForecast.Range("H6").Value = 8.56
$wsForecast.Range("L6").Value = 0.85

# Row 7
$wsForecast.Range("H7").Value = 7.97
$wsForecast.Range("L7").Value = 1.01

# Row 8
$wsForecast.Range("H8").Value = 6.88

# Row 9
$wsForecast.Range("H9").Value = 5.88
$wsForecast.Range("L9").Value = 1.18

# Row 10
$wsForecast.Range("D10").Value = 87
$wsForecast.Range("H10").Value = 5.54

# Row 11
$wsForecast.Range("D11").Value = 92
$wsForecast.Range("H11").Value = 4.29
$wsForecast.Range("L11").Value = 0.96

# Row 12
$wsForecast.Range("H12").Value = 3.23
$wsForecast.Range("L12").Value = 0.99

# Row 13
$wsForecast.Range("H13").Value = 2.17
$wsForecast.Range("L13").Value = 0.9

# Row 14
$wsForecast.Range("H14").Value = 1.26
$wsForecast.Range("J14").Value = "Normal"
$wsForecast.Range("L14").Value = 0.95

# Row 15
$wsForecast.Range("D15").Value = 78
$wsForecast.Range("H15").Value = 0.3
$wsForecast.Range("L15").Value = 0.92

# Row 16
$wsForecast.Range("D16").Value = 69
$wsForecast.Range("L16").Value = 0.99

# Row 17
$wsForecast.Range("D17").Value = 71
$wsForecast.Range("L17").Value = 0.82

# --- Summary sheet updates ---
# These "Value" cells are stored as text (numeric-looking strings), so a
# leading apostrophe forces text entry instead of Excel auto-converting to
# a number (mirrors how the source file stores them as inline strings).
# Re-apply the "Normal" style afterward so the quote-prefix formatting
# Excel applies for the apostrophe doesn't leave a stray cell style behind.

$wsSummary.Range("B9").Value = "'1479"
$wsSummary.Range("B9").Style = "Normal"

$wsSummary.Range("B10").Value = "'802"
$wsSummary.Range("B10").Style = "Normal"

$wsSummary.Range("B11").Value = "'404"
$wsSummary.Range("B11").Style = "Normal"

$wsSummary.Range("B12").Value = "'124"
$wsSummary.Range("B12").Style = "Normal"

$wsSummary.Range("B14").Value = "'69"
$wsSummary.Range("B14").Style = "Normal"
